$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain their original text formatting
# (values such as "0.9998" or "0.07650" must not be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.238.55'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.842.49'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '242.24'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").Value = '0.6639'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.9994'
$ws.Range("D8").Value = '0.07457'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").Value = '0.2952'
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").Value = '23.36'
$ws.Range("E10").Value = '  +2.72%  '
$ws.Range("D11").Value = '0.07763'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").Value = '1.846.61'
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").Value = '5.025'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").Value = '0.6735'
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").Value = '83.40'
$ws.Range("E15").Value = '  -3.06%  '
$ws.Range("D16").Value = '6.185'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '0.000008621'
$ws.Range("E17").Value = '  +4.84%  '
$ws.Range("D18").Value = '29.239.48'
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").Value = '2.096.78'
$ws.Range("E19").Value = '  +1.43%  '
$ws.Range("D20").Value = '228.49'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").Value = '7.191'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").Value = '0.9999'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '159.10'
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("D26").Value = '0.1413'
$ws.Range("E26").Value = '  +1.39%  '
$ws.Range("D27").Value = '8.638'
$ws.Range("D28").Value = '18.08'
$ws.Range("D29").Value = '1.511'
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("D30").Value = '4.133'
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("D31").Value = '4.058'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("D33").Value = '0.05337'
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("D34").Value = '1.885'
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("D35").Value = '0.7450'
$ws.Range("E35").Value = '  -0.56%  '
$ws.Range("D36").Value = '1.156'
$ws.Range("E36").Value = '  +2.17%  '
$ws.Range("D37").Value = '2.656'
$ws.Range("E37").Value = '  -0.85%  '
$ws.Range("D38").Value = '1.323.14'
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("D39").Value = '0.01798'
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("D40").Value = '2.745'
$ws.Range("E40").Value = '  +1.00%  '
$ws.Range("D41").Value = '6.433'
$ws.Range("E41").Value = '  +7.73%  '
$ws.Range("D42").Value = '0.9205'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '0.9990'
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = '103.20'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").Value = '66.29'
$ws.Range("E45").Value = '  +4.02%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.995.93'
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.00000000124'
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.5139'
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("B49").Value = 'XinFinNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D49").Value = '0.07650'
$ws.Range("E49").Value = '  -5.48%  '
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("E51").Value = '  -1.24%  '
